$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the used range (content + formatting),
# then rebuild the table with its new columns/values/formats.
$ws.UsedRange.Clear()

# ---- Header row (bold) ----
$ws.Range("A1").Value = "Level"
$ws.Range("B1").Value = "Class"
$ws.Range("C1").Value = "Labeling error % (Felipe)"
$ws.Range("D1").Value = "Labeling error % (Karen)"
$ws.Range("E1").Value = "Labeling error % (Avg)"
$ws.Range("F1").Value = "F1 Score"
$ws.Range("G1").Value = "Positive Auto-Manual % diff"
$ws.Range("H1").Value = "Correlation F1-Labeling Error"
$ws.Range("A1:H1").Font.Bold = $true

$ws.Range("H2").Value = "Correlation F1-Positive Diff"
$ws.Range("H2").Font.Bold = $true

# ---- Correlation formulas (column I) ----
$ws.Range("I1").Formula = "=CORREL(E2:E8,F2:F8)"
$ws.Range("I1").NumberFormat = "0.00"
$ws.Range("I2").Formula = "=CORREL(F2:F8,G2:G8)"
$ws.Range("I2").NumberFormat = "0.00"
$ws.Range("I3:I8").NumberFormat = "0.00"

# ---- Level numbers (A2:A8), matching the original merge A4:A8 ----
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A2:A3").HorizontalAlignment = -4108
$ws.Range("A2:A3").VerticalAlignment = -4108
$ws.Range("A4:A8").HorizontalAlignment = -4108
$ws.Range("A4:A8").VerticalAlignment = -4108
$ws.Range("A4:A8").Merge()

# ---- Event names (Class column, B) ----
$ws.Range("B2").Value = "Duct"
$ws.Range("B3").Value = "Event"
$ws.Range("B4").Value = "Anode"
$ws.Range("B5").Value = "Buried"
$ws.Range("B6").Value = "Damage"
$ws.Range("B7").Value = "Flange"
$ws.Range("B8").Value = "Repair"

# ---- Labeling error % (Felipe) - column C ----
$ws.Range("C2").Value = 0.0757
$ws.Range("C3").Value = 0.1006
$ws.Range("C4").Value = 0.0244
$ws.Range("C5").Value = 0.1452
$ws.Range("C6").Value = 0.0534
$ws.Range("C7").Value = 0.1185
$ws.Range("C8").Value = 0.9511

# ---- Labeling error % (Karen) - column D ----
$ws.Range("D2").Value = 0.0757
$ws.Range("D3").Value = 0.1205
$ws.Range("D4").Value = 0.0375
$ws.Range("D5").Value = 0.2102
$ws.Range("D6").Value = 0.1386
$ws.Range("D7").Value = 0.0688
$ws.Range("D8").Value = 0.9453

# ---- Labeling error % (Avg) - column E ----
$ws.Range("E2").Formula = "=AVERAGE(C2:D2)"
$ws.Range("E3").Formula = "=AVERAGE(C3:D3)"
$ws.Range("E4").Formula = "=AVERAGE(C4:D4)"
$ws.Range("E5").Formula = "=AVERAGE(C5:D5)"
$ws.Range("E6").Formula = "=AVERAGE(C6:D6)"
$ws.Range("E7").Formula = "=AVERAGE(C7:D7)"
$ws.Range("E8").Formula = "=AVERAGE(C8:D8)"

$ws.Range("C2:E9").NumberFormat = "0.00%"

# ---- F1 Score - column F ----
$ws.Range("F2").Value = 0.93598862019914653
$ws.Range("F3").Value = 0.87889273356401376
$ws.Range("F4").Value = 0.69699999999999995
$ws.Range("F5").Value = 0.41299999999999998
$ws.Range("F6").Value = 0.81100000000000005
$ws.Range("F7").Value = 0.72899999999999998
$ws.Range("F8").Value = 0.247
$ws.Range("F2:F8").NumberFormat = "0.00"

# ---- Positive Auto-Manual % diff - column G ----
$ws.Range("G2").Value = 0.08
$ws.Range("G3").Value = 0.17
$ws.Range("G4").Value = 0.030765086669152203
$ws.Range("G5").Value = 0.32173350150617291
$ws.Range("G6").Value = 0.26669531662415558
$ws.Range("G7").Value = 0.24805868053601665
$ws.Range("G8").Value = 0.74578678005320553
$ws.Range("G2:G8").NumberFormat = "0%"

# ---- Footer label ----
$ws.Range("A9").Value = "Semiauto dataset"

# ---- Leftover formatted-but-empty placeholder cells (column H) ----
$ws.Range("H11").NumberFormat = "0%"
$ws.Range("H12").NumberFormat = "0%"
$ws.Range("H14").NumberFormat = "0%"
$ws.Range("H15").NumberFormat = "0%"
$ws.Range("H17").NumberFormat = "0%"
$ws.Range("H18").NumberFormat = "0%"
$ws.Range("H20").NumberFormat = "0%"
$ws.Range("H21").NumberFormat = "0%"
$ws.Range("H23").NumberFormat = "0%"
$ws.Range("H24").NumberFormat = "0%"

# ---- Column widths to fit the new headers ----
$ws.Columns.Item(2).ColumnWidth = 12.1796875
$ws.Columns.Item(3).ColumnWidth = 21.1796875
$ws.Columns.Item(4).ColumnWidth = 21.08984375
$ws.Columns.Item(5).ColumnWidth = 19.1796875
$ws.Columns.Item(7).ColumnWidth = 24.54296875
$ws.Columns.Item(8).ColumnWidth = 25

# ---- Threaded comment from Olavo Sampaio on B2 ----
$excel.UserName = "Olavo Sampaio"
$ws.Range("B2").AddCommentThreaded("anotado apenas pelo Felipe")

$ws.Range("H6").Select()
